# Reorder the "Recorded By" (column G) comma-separated contributor lists so
# that any exact "System" entry is moved to the end of the list, preserving
# the relative order of the remaining entries (including a lowercase
# "system" entry, which is left in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -ne "") {
        $parts = $text.Split(",")

        if ($parts.Length -gt 1) {
            $others = @()
            $systems = @()

            foreach ($p in $parts) {
                $trimmed = $p.Trim()
                if ($trimmed.Equals("System")) {
                    $systems += $trimmed
                } else {
                    $others += $trimmed
                }
            }

            if ($systems.Length -gt 0 -and $others.Length -gt 0) {
                $newParts = $others + $systems
                $newText = [string]::Join(", ", $newParts)
                if (-not $newText.Equals($text)) {
                    $cell.Value = $newText
                }
            }
        }
    }
}
